# Daily attendance processing - 2025-12-01 22:24:54
# Swap the order of the two comma-separated names/emails in the
# "Recorded By" column (G) so that "dnasr281@gmail.com" is listed first,
# for every row where it currently appears second (e.g. "System, dnasr281@gmail.com"
# -> "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(", ")) {
        $parts = $val.Split(", ")
        if ($parts.Length -eq 2 -and $parts[1] -eq "dnasr281@gmail.com" -and $parts[0] -ne "dnasr281@gmail.com") {
            $cell.Value2 = "dnasr281@gmail.com, " + $parts[0]
        }
    }
}
